# Edit script: add two new test-case rows (rows 18-21) to the "Konto" sheet,
# mirroring the formatting of the existing "Poprawne wylogowanie" block
# (rows 16-17), and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Block 1: "Dodawanie produktu do koszyka" (rows 18-19) ---------------

# Copy the formatting of the previous test-case block (rows 16-17) down
# onto the two new rows, then fill in the block-specific formatting tweaks
# (extra wrap-text) before writing values, so Excel doesn't re-wrap/resize
# things out from under us.
$ws.Range("A16:H17").Copy()
$ws.Range("A18:H19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B18:B19").WrapText = $true
$ws.Range("G18:G19").WrapText = $true
$ws.Range("H18:H19").WrapText = $true

$ws.Range("A18:A19").Merge()
$ws.Range("B18:B19").Merge()
$ws.Range("C18:C19").Merge()
$ws.Range("F18:F19").Merge()
$ws.Range("G18:G19").Merge()
$ws.Range("H18:H19").Merge()

$ws.Cells.Item(18, 2).Value = "Dodawanie produktu do koszyka"
$ws.Cells.Item(18, 3).Value = "Zalogowany użytkownik z użyciem danych: 235689klasa@gmail.com;  A1@a2222"
$ws.Cells.Item(18, 4).Value = '1. Ze strony głównej dodaj do koszyka  "plusem" dostepny, losowy produkt w ilości 1.'
$ws.Cells.Item(18, 5).Value = "Cosma Original Snackies"
$ws.Cells.Item(18, 6).Value = "Wybrany produkt zostaje dodany do koszyka zakupowego. "
$ws.Cells.Item(18, 7).Value = "Produkt znajduję się w koszyku."
$ws.Cells.Item(18, 8).Value = "Pass"

$ws.Cells.Item(19, 4).Value = "2. Przejdź do koszyka i zweryfikuje, czy produkt się w nim znajduje."

$ws.Rows.Item(18).RowHeight = 15
$ws.Rows.Item(19).RowHeight = 50.25

# --- Block 2: "Wyszukiwanie produktu" (rows 20-21) ------------------------

$ws.Range("A16:H17").Copy()
$ws.Range("A20:H21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A20:A21").Merge()
$ws.Range("B20:B21").Merge()
$ws.Range("C20:C21").Merge()
$ws.Range("F20:F21").Merge()
$ws.Range("G20:G21").Merge()
$ws.Range("H20:H21").Merge()

$ws.Cells.Item(20, 2).Value = "Wyszukiwanie produktu "
$ws.Cells.Item(20, 3).Value = "Zalogowany użytkownik z użyciem danych: 235689klasa@gmail.com;  A1@a2222"
$ws.Cells.Item(20, 4).Value = "1. Wpisz w pole wyszukiwania słowo kangur"
$ws.Cells.Item(20, 5).Value = "kangur"
$ws.Cells.Item(20, 6).Value = "System zwraca produkty, które mają wyszukiwane słowo w opisie produtu. "
$ws.Cells.Item(20, 7).Value = 'Wyświetlone produkty mają słowo "kangur" w opisie. '
$ws.Cells.Item(20, 8).Value = "Pass"

$ws.Cells.Item(21, 4).Value = '2. Zweryfikuj, czy wyświetlone produkty zawierają słowo "kangur" w opisie.'

$ws.Rows.Item(20).RowHeight = 15
$ws.Rows.Item(21).RowHeight = 15.75

# --- Selection -------------------------------------------------------------

$ws.Range("C16:C17").Select()
